$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: after "source1_cleaned.csv" insert " (A.csv in the Data folder)"
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("source1_cleaned.csv", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Collapse(0)
    $r.InsertAfter(" (A.csv in the Data folder)")
}

# ---------------------------------------------------------------------------
# Change 2: after "source2_cleaned.csv" insert " (B.csv in the Data folder)"
# plus a _GoBack bookmark right after it (moved from the hyperlink below)
# ---------------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("source2_cleaned.csv", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $r2.Collapse(0)
    $r2.InsertAfter(" (B.csv in the Data folder)")
    $r2.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $r2) | Out-Null
}

# ---------------------------------------------------------------------------
# Change 3: fix the "here" hyperlink (previously split across runs "h"/"e"/"re"
# with a _GoBack bookmark in between) into one clean run, no bookmark.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks.Item($i)
    if ($h.TextToDisplay -eq "here") {
        $h.TextToDisplay = "here"
        break
    }
}

# ---------------------------------------------------------------------------
# Change 4: append a new sentence after "6.5 seconds"
# ---------------------------------------------------------------------------
$r4 = $d.Content
$found4 = $r4.Find.Execute("6.5 seconds", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
    $r4.Collapse(0)
    $r4.InsertAfter(". We did, however, try out different blockers, analysed the output and loosened/tightened them. The total time for that took ~2 hours. (The entire process from running, debugging and iterating)")
}

# ---------------------------------------------------------------------------
# Change 5: rewrite the "Since we did not have to debug..." bullet, and set
# bold/size-24 formatting on the paragraph mark (w:pPr/w:rPr) only.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Since we did not have to debug*") {
        $xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Each step of debugging took around 1 hour  (running the debugger, looking at false positives/negatives, relabeling if required/trying out different feature combinations/adding triggers). We had 3 main iterations, so the total time was 3 hours. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
        $p.Range.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# Change 6: add a lastRenderedPageBreak before "For the test set, we obtained..."
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "For the test set, we obtained*") {
        $xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>For the test set, we obtained perfect recall. This doesn&#8217;t say enough about the matcher being perfect because &#8211;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
        $p.Range.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# Change 7: remove the lastRenderedPageBreak from "For the training set recall..."
# (keep the second run / sentence that follows in the same paragraph)
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "For the training set recall*") {
        $xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">For the training set recall, we believe that we couldn&#8217;t increase it more because of difficult false negatives. </w:t></w:r><w:r><w:t>For example, there were some book pairs which were the same, but in the ltable, their name was truncated, but the rtable, the name for more descriptive. For example &#8211;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
        $p.Range.InsertXML($xml)
        break
    }
}

Write-Host "All changes applied"
